# Updated a few README.md & Removed the scripts/models/ directory.
# It was useless.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sign flips on existing "Ecart" (H/I) columns, rows 2-5 -----------------
$ws.Range("H2").Value = -135
$ws.Range("I3").Value = -15
$ws.Range("I4").Value = -15
$ws.Range("H5").Value = -100

# --- New rows of servo data --------------------------------------------------
# Row 6: AVMG
$ws.Range("A6").Value = "AVMG"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 220
$ws.Range("F6").Value = 90
$ws.Range("G6").Value = 400
$ws.Range("H6").Value = 180
$ws.Range("I6").Value = -130

# Row 8: ARMG (written before row 7's AVMD so shared-string order matches)
$ws.Range("A8").Value = "ARMG"
$ws.Range("B8").Value = 14
$ws.Range("C8").Value = 340
$ws.Range("F8").Value = 200
$ws.Range("G8").Value = 505
$ws.Range("H8").Value = 165
$ws.Range("I8").Value = -140

# Row 7: AVMD
$ws.Range("A7").Value = "AVMD"
$ws.Range("B7").Value = 14
$ws.Range("C7").Value = 460
$ws.Range("G7").Value = 290
$ws.Range("H7").Value = -170
$ws.Range("I7").Value = 95

# Row 13 / 14 labels (entered here to match original authoring order)
$ws.Range("D13").Value = "Min: Le plus haut"
$ws.Range("D14").Value = "Max: Le plus bas"

# Back to row 7's F cell (text value instead of a number)
$ws.Range("F7").Value = "555 (angle non parfait)"

# Row 9: ARMD / servo replacement note
$ws.Range("A9").Value = "ARMD"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "REPLACER LE SERVO (75)"

# Row 10: AVJG
$ws.Range("A10").Value = "AVJG"
$ws.Range("B10").Value = 0

# --- Selection, matching the saved workbook state ----------------------------
$ws.Range("K38").Select() | Out-Null
